$d = $word.ActiveDocument

# 1. "ima u razvoju projekta i testiranju..." -> remove the space so it reads
#    "ima u razvojuprojekta i testiranju..."
$d.Content.Find.Execute("ima u razvoju projekta i testiranju", $true, $false, $false, $false, $false, $true, 1, $false, "ima u razvojuprojekta i testiranju", 2)

# 2. "...u koraku 4 bira tab sa natpisom "Najteže" i dobija pregled..." -> remove
#    the two spaces so it reads "...u koraku 4bira tab sa natpisom"Najteže" i dobija pregled..."
$d.Content.Find.Execute("u koraku 4 bira tab sa natpisom ", $true, $false, $false, $false, $false, $true, 1, $false, "u koraku 4bira tab sa natpisom", 2)

# 3. "...ime autora čije ideje želi..." -> remove the space so it reads
#    "...ime autora čijeideje želi..."
$d.Content.Find.Execute("čije ideje želi", $true, $false, $false, $false, $false, $true, 1, $false, "čijeideje želi", 2)
